$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = 112196385
$ws.Range("B8").Value = 89423
$ws.Range("C8").Value = "Ovaliderad"
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 5432
$ws.Range("F8").Value = "Granticka"
$ws.Range("G8").Value = "Porodaedalea chrysoloma"
$ws.Range("H8").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""
$ws.Range("K8").Value = ""
$ws.Range("N8").Value = ""
$ws.Range("P8").Value = "Jeppmakullarna, Ly lm"
$ws.Range("Q8").Value = 606897.2701169839
$ws.Range("R8").Value = 7189526.380129344
$ws.Range("S8").Value = 25
$ws.Range("T8").Value = "Västerbotten"
$ws.Range("U8").Value = "Storuman"
$ws.Range("V8").Value = "Lycksele lappmark"
$ws.Range("W8").Value = "Stensele"
$ws.Range("Y8").Value = "'2023-08-17"
$ws.Range("Z8").Value = "00:00"
$ws.Range("AA8").Value = "'2023-08-17"
$ws.Range("AB8").Value = "00:00"
$ws.Range("AD8").Value = $false
$ws.Range("AE8").Value = $false
$ws.Range("AF8").Value = ""
$ws.Range("AG8").Value = $false
$ws.Range("AT8").Value = ""
$ws.Range("AW8").Value = "Lena Lundevaller"
$ws.Range("AX8").Value = "Lena Lundevaller, Åsa Stenman"
$ws.Range("AY8").Value = ""
